$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Text updates: explicitly call out the programming language used
#    alongside each framework/tool mention.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Python, Java, JavaScript (Node.js), C#, Bash, CSS, HTML",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Python, Java, JavaScript, Node.js, C#, Bash, CSS, HTML", 2)

$d.Content.Find.Execute(
    "Enhanced legacy ASP.NET codebase to enable political campaigns to organize volunteers and fundraise.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Enhanced legacy C#/ASP.NET codebase to enable political campaigns to organize volunteers and fundraise.", 2)

$d.Content.Find.Execute(
    "Designed and built React/Redux applications to interface with HTTP APIs.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Designed and built JavaScript/React applications to interface with HTTP APIs.", 2)

$d.Content.Find.Execute(
    "Deployed a Flask API for internal use on a Gunicorn server and Nginx reverse proxy.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Deployed a Python/Flask API for internal use on a Gunicorn server and Nginx reverse proxy.", 2)

$d.Content.Find.Execute(
    "Automated portions of the event coordination process by writing a Selenium-based web crawler.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Automated portions of the event coordination process by writing a Python/Selenium web scraper.", 2)

# ---------------------------------------------------------------------
# 2. Page margins: left/right reduced from 0.6" (864 twips) to
#    0.5" (720 twips). This widens the usable body width from 10512
#    to 10800 twips, which is why every full-width table below grows.
# ---------------------------------------------------------------------
$d.Sections(1).PageSetup.LeftMargin = 36
$d.Sections(1).PageSetup.RightMargin = 36

# ---------------------------------------------------------------------
# 3. Tables: every table spans 100% of the now-wider body, so each
#    column grows proportionally. Widths below are in points
#    (twips / 20) to match the target dxa values from the edit.
# ---------------------------------------------------------------------

# Table 1 - EDUCATION: University of Connecticut
$t1 = $d.Tables(1)
$t1.Columns(1).Width = 361.3
$t1.Columns(2).Width = 178.65

# Table 2 - EDUCATION: Harvard University
$t2 = $d.Tables(2)
$t2.Columns(1).Width = 361.25
$t2.Columns(2).Width = 178.7

# Table 3 - TECHNICAL SKILLS
$t3 = $d.Tables(3)
$t3.Columns(1).Width = 129.75
$t3.Columns(2).Width = 410.2

# Table 4 - EXPERIENCE: Software Engineering Intern (Full-Time)
$t4 = $d.Tables(4)
$t4.Columns(1).Width = 361.3
$t4.Columns(2).Width = 178.65

# Table 5 - EXPERIENCE: Software Engineer
$t5 = $d.Tables(5)
$t5.Columns(1).Width = 361.3
$t5.Columns(2).Width = 178.65

# Table 6 - EXPERIENCE: Full Stack Developer
$t6 = $d.Tables(6)
$t6.Columns(1).Width = 361.3
$t6.Columns(2).Width = 178.65
